# Version 0.2.01 (Alpha) - Two Indicators & Only UP (Rate 1).
#
# The sheet holds two mirrored "indicator" tables (D:E and J:K). This
# edit removes the empty spacer column (I) between them - which shifts
# the second table from J:K to I:J - and appends a new "Captured Val: "
# label row (19) to both tables, copying the formatting of the last
# existing row (18).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the empty spacer column I - this shifts the right-hand
#    table's columns J->I and K->J (values, styles, merged cell D1:E1
#    style J1:K1 -> I1:J1, etc. all move automatically).
$ws.Columns.Item(9).Delete()

# 2) Append row 19 to both tables, copying row 18's cell formatting.
$ws.Range("D18:E18").Copy()
$ws.Range("D19:E19").PasteSpecial(-4122)

$ws.Range("I18:J18").Copy()
$ws.Range("I19:J19").PasteSpecial(-4122)

$ws.Range("D19").Value = "Captured Val: "
$ws.Range("I19").Value = "Captured Val: "

# 3) Match the resized column widths from the author's edit.
$ws.Columns.Item(4).ColumnWidth = 12.022135416666666
$ws.Columns.Item(9).ColumnWidth = 12.166666666666666
$ws.Columns.Item(10).ColumnWidth = 10.166666666666666
$ws.Columns.Item(11).ColumnWidth = 9.022135416666666

# 4) Leave the selection where the author left theirs.
[void]$ws.Range("K22").Select()
